# Generate Report for Handoff
#
# A new handoff run produced a freshly generated GUID-named markdown /
# xliff pair. Replace every occurrence of the old GUID-based file names
# and refresh the handoff / handback timestamps that the report records,
# across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "b69c3bdf-1046-4d6e-b7e3-c96289c64ffe"
$newGuid = "74ea9e28-b2f5-4609-a8c7-113c66e10282"

$newHoDate = "2016-08-13 05:10:19"

$newZhHash = "03ea18f8c25e25beffb2065396c6f3a6e7f93e22"
$newZhDate = "2016-08-13 05:10:11"

$newDeHash = "03ea18f8c25e25beffb2065396c6f3a6e7f93e22"

# The external hyperlink target (the GitHub blob URL) itself is not part
# of this edit - only the displayed text changes - so reuse the address
# that is already on the sheet's single hyperlink.
$linkAddress = "https://github.com/OpenLocalizationTestOrg/oltest/blob/4faeb37fa2ede4379e4e403f6c08bf497e1b4946/e2e/$oldGuid.md"

function Update-Hyperlink {
    param($ws, $cellAddr, $display)

    $cell = $ws.Range($cellAddr)

    # This host always appends a fresh <hyperlink> entry when a property
    # of an existing Hyperlinks.Item(...) is written (it never mutates in
    # place), so clear the cell's hyperlink(s) first and re-add a single
    # clean one with the updated display text.
    $cell.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($cell, $linkAddress, "", "", $display)

    # Hyperlinks.Add re-stamps the cell with the host's built-in
    # "Hyperlink" style; restore the workbook's original look (underlined,
    # cornflower-blue text) so no unrelated styling changes leak in.
    # Font.Color takes a BGR-packed value (OLE_COLOR order), so RGB
    # 6495ED ("cornflower blue") is written as 0xED9564.
    $cell.Font.Underline = 2
    $cell.Font.Color = 0xED9564
}

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = $newHoDate

Update-Hyperlink $wsOverview "B2" "e2e\$newGuid.md"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newZhHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = $newZhDate

Update-Hyperlink $wsZhCn "A2" "$newGuid.md"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newDeHash.de-de.xlf"
$wsDeDe.Range("H2").Value = $newHoDate

Update-Hyperlink $wsDeDe "A2" "$newGuid.md"
